$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3, shifting existing rows (3..51) down to (4..52)
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new claim data.
# Columns A-L are stored as text; force text number format so values
# (including ones that look numeric/date-like) are written as strings,
# then reset the style back to Normal so no extra formatting is applied.
$textRange = $ws.Range("A3:L3")
$textRange.NumberFormat = "@"

$ws.Range("A3").Value = "1939"
$ws.Range("B3").Value = "5/8/2024"
$ws.Range("C3").Value = "SANCHEZ, MIGUEL B. AV. 1050"
$ws.Range("D3").Value = "13"
$ws.Range("E3").Value = "788825789"
$ws.Range("F3").Value = "Optical Power"
$ws.Range("G3").Value = "Pendiente"
$ws.Range("H3").Value = "Dos postes rajados y una columna fuera de plomo. Ver fotos o pedirme ubicacion, esta en la puerta de un colegio"
$ws.Range("I3").Value = "1"
$ws.Range("J3").Value = "Cambio"
$ws.Range("K3").Value = "Sin equipos"
$ws.Range("L3").Value = "Poste"

$textRange.Style = "Normal"

# Columns M and N are numeric coordinates.
$ws.Range("M3").Value = -58.455394
$ws.Range("N3").Value = -34.542575
